$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.697.86"
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = "'1.728.18"
$ws.Range("E3").Value = '  +0.06%  '
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = "'242.00"
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("D6").Value = "'0.9979"
$ws.Range("E6").Value = '  -0.29%  '
$ws.Range("D7").Value = "'0.4922"
$ws.Range("E7").Value = '  +0.65%  '
$ws.Range("D8").Value = "'0.2618"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'0.06219"
$ws.Range("E9").Value = '  +0.43%  '
$ws.Range("D10").Value = "'1.723.81"
$ws.Range("E10").Value = '  -0.20%  '
$ws.Range("D11").Value = "'15.81"
$ws.Range("E11").Value = '  +2.21%  '
$ws.Range("D12").Value = "'0.06986"
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = "'0.6102"
$ws.Range("E13").Value = '  +1.74%  '
$ws.Range("D14").Value = "'4.491"
$ws.Range("E14").Value = '  -1.73%  '
$ws.Range("D15").Value = "'77.32"
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("D17").Value = "'26.505.73"
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = "'0.9978"
$ws.Range("E18").Value = '  -0.30%  '
$ws.Range("D19").Value = "'0.000007238"
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("D20").Value = "'11.39"
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = "'1.950.93"
$ws.Range("E21").Value = '  -0.30%  '
$ws.Range("D22").Value = "'4.456"
$ws.Range("E22").Value = '  -0.73%  '
$ws.Range("D23").Value = "'8.550"
$ws.Range("E23").Value = '  -0.69%  '
$ws.Range("D24").Value = "'5.085"
$ws.Range("E24").Value = '  -2.12%  '
$ws.Range("D25").Value = "'138.02"
$ws.Range("E25").Value = '  -0.70%  '
$ws.Range("D26").Value = "'15.33"
$ws.Range("E26").Value = '  +0.30%  '
$ws.Range("D27").Value = "'1.773"
$ws.Range("E27").Value = '  +2.83%  '
$ws.Range("D28").Value = "'106.43"
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("E29").Value = '  -1.52%  '
$ws.Range("D30").Value = "'3.921"
$ws.Range("E30").Value = '  -1.12%  '
$ws.Range("D31").Value = "'0.07974"
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("D32").Value = "'3.669"
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").Value = "'0.04465"
$ws.Range("E33").Value = '  -1.34%  '
$ws.Range("D34").Value = "'0.9974"
$ws.Range("D35").Value = "'2.610"
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("D37").Value = "'0.6235"
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("D38").Value = "'0.9319"
$ws.Range("E38").Value = '  +2.41%  '
$ws.Range("D39").Value = "'2.044"
$ws.Range("E39").Value = '  +3.04%  '
$ws.Range("D40").Value = "'2.408"
$ws.Range("E40").Value = '  +0.25%  '
$ws.Range("D41").Value = "'0.9978"
$ws.Range("E41").Value = '  -0.37%  '
$ws.Range("D42").Value = "'0.01512"
$ws.Range("E42").Value = '  +1.61%  '
$ws.Range("D43").Value = "'5.610"
$ws.Range("E43").Value = '  +3.60%  '
$ws.Range("D44").Value = "'99.36"
$ws.Range("E44").Value = '  -1.01%  '
$ws.Range("D45").Value = "'0.3847"
$ws.Range("E45").Value = '  -0.57%  '
$ws.Range("D46").Value = "'6.870"
$ws.Range("E46").Value = '  +2.44%  '
$ws.Range("D47").Value = "'0.1157"
$ws.Range("E47").Value = '  -0.12%  '
$ws.Range("D48").Value = "'0.05380"
$ws.Range("E48").Value = '  +0.29%  '
$ws.Range("D49").Value = "'7.856"
$ws.Range("E49").Value = '  +2.41%  '
$ws.Range("D50").Value = "'30.24"
$ws.Range("E50").Value = '  -0.28%  '
$ws.Range("D51").Value = "'51.70"
$ws.Range("E51").Value = '  +1.17%  '
